# Natmi following Dr Hou advice
# Update Il12a-Il12rb1 LR-pairs sheet: refresh row 2 (sCs -> ECs) with new
# stats, and add three more target-cluster rows (FAPs, M1, M2) plus a
# self-referencing sCs row, each carrying the refreshed ligand/edge stats.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: sCs | Il12a | Il12rb1 | ECs
$ws.Range("A2").Value = "sCs"
$ws.Range("B2").Value = "Il12a"
$ws.Range("C2").Value = "Il12rb1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.8401833333333332
$ws.Range("H2").Value = 2.52055
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.1866496666666667
$ws.Range("N2").Value = 0.559949
$ws.Range("O2").Value = 0.1052515138801798
$ws.Range("P2").Value = 0.1052515138801798
$ws.Range("Q2").Value = 0.1568199391055555
$ws.Range("R2").Value = 1.41137945195
$ws.Range("S2").Value = 0.1052515138801798
$ws.Range("T2").Value = 0.1052515138801798

# Row 3: sCs | Il12a | Il12rb1 | FAPs
$ws.Range("A3").Value = "sCs"
$ws.Range("B3").Value = "Il12a"
$ws.Range("C3").Value = "Il12rb1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.8401833333333332
$ws.Range("H3").Value = 2.52055
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.9694063333333333
$ws.Range("N3").Value = 2.908219
$ws.Range("O3").Value = 0.5466470204341869
$ws.Range("P3").Value = 0.546647020434187
$ws.Range("Q3").Value = 0.8144790444944442
$ws.Range("R3").Value = 7.330311400449999
$ws.Range("S3").Value = 0.5466470204341869
$ws.Range("T3").Value = 0.546647020434187

# Row 4: sCs | Il12a | Il12rb1 | M1
$ws.Range("A4").Value = "sCs"
$ws.Range("B4").Value = "Il12a"
$ws.Range("C4").Value = "Il12rb1"
$ws.Range("D4").Value = "M1"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.8401833333333332
$ws.Range("H4").Value = 2.52055
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.112113
$ws.Range("N4").Value = 0.336339
$ws.Range("O4").Value = 0.06322038065421277
$ws.Range("P4").Value = 0.06322038065421277
$ws.Range("Q4").Value = 0.09419547404999999
$ws.Range("R4").Value = 0.8477592664499999
$ws.Range("S4").Value = 0.06322038065421277
$ws.Range("T4").Value = 0.06322038065421277

# Row 5: sCs | Il12a | Il12rb1 | M2
$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Il12a"
$ws.Range("C5").Value = "Il12rb1"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.8401833333333332
$ws.Range("H5").Value = 2.52055
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.30418
$ws.Range("N5").Value = 0.91254
$ws.Range("O5").Value = 0.1715267220340053
$ws.Range("P5").Value = 0.1715267220340054
$ws.Range("Q5").Value = 0.2555669663333333
$ws.Range("R5").Value = 2.300102697
$ws.Range("S5").Value = 0.1715267220340053
$ws.Range("T5").Value = 0.1715267220340054

# Row 6: sCs | Il12a | Il12rb1 | sCs
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Il12a"
$ws.Range("C6").Value = "Il12rb1"
$ws.Range("D6").Value = "sCs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.8401833333333332
$ws.Range("H6").Value = 2.52055
$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 1
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.201019
$ws.Range("N6").Value = 0.603057
$ws.Range("O6").Value = 0.1133543629974151
$ws.Range("P6").Value = 0.1133543629974151
$ws.Range("Q6").Value = 0.1688928134833333
$ws.Range("R6").Value = 1.52003532135
$ws.Range("S6").Value = 0.1133543629974151
$ws.Range("T6").Value = 0.1133543629974151
